$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("m05")

# Insert the two new rows first (shifts existing rows down).
[void]$ws.Rows.Item(2).Insert()
[void]$ws.Rows.Item(36).Insert()

# Fill column A and C (mirrored name) for both new rows first, matching the
# order new shared strings were originally introduced.
$ws.Range("A2").Value2 = "DX_M05_0005_trent"
$ws.Range("C2").Value2 = "DX_M05_0005_trent"
$ws.Range("A36").Value2 = "DX_M05_0395_trent"
$ws.Range("C36").Value2 = "DX_M05_0395_trent"

# Then column B (hash codes) for the new rows, plus the single updated hash
# on the pre-existing "DX_M05_0150_smuggler2" row (now row 16).
$ws.Range("B36").Value2 = "0x8B1F4040"
$ws.Range("B2").Value2 = "0xAB426743"
$ws.Range("B16").Value2 = "0xA7C4FF08"

# Keep the hidden filter-database defined name in sync with the new extent.
$wb.Names.Item("m05!_FilterDatabase").RefersTo = "='m05'!`$A`$1:`$C`$94"

# Reproduce the author's final on-sheet selection.
[void]$ws.Activate()
[void]$ws.Range("B11:B16").Select()
